# Generate Report for handoff
#
# The "41cbefcb-0561-4345-abaf-f0b90f69efd9" file has finished its
# handback cycle and is dropped from the report; the "0d410992..." file
# moves on to "Ready for handoff" with refreshed handback timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the 41cbefcb-...md row (row 3) -- .localization-config shifts up
# from row 4 to row 3, carrying its original formatting with it.
$ws1.Rows.Item(3).Delete()

# 0d410992-...md is now ready for handoff again.
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Hyperlinks aren't renumbered by the row delete, so rebuild the set.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Delete()

# Refresh the Latest Handback DateTime for 0d410992-...
$ws2.Range("D2").Value = "2016-01-13 04:48:14"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20e98bc58e59e0e4d9414e6326a3d2ea724a7362/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9901b54e415d3d10064531cddf6c6cb378273d2d/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/02146fa5a99fc537bdcf65c69e38fb70436f77b1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Delete()

# Refresh the Latest Handback DateTime for 0d410992-...
$ws3.Range("D2").Value = "2016-01-13 04:48:35"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d527aa0fd447ebb9602d643ba23e6a7e6362ea59/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/571f56c13ca65932256543eee4dd4566f0a865b1/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/78b219b4b8eb6334c8c716cf0fef493ef5875863/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf", "", "", "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed60080193a1e35195b955bf39917ca9abd1f2cf/.localization-config", "", "", ".localization-config") | Out-Null

"Report regenerated for handoff"
